$d = $word.ActiveDocument

# =========================================================================
# 1. Bold "We regularly collaborate" inside the "OPENING TEXT" paragraph,
#    leaving the remainder of the sentence (" with other functions ...")
#    exactly as it was.
# =========================================================================

# Locate the paragraph dynamically (don't hard-code an index).
$targetIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "We regularly collaborate*") {
        $targetIdx = $i
        break
    }
}

if ($targetIdx -gt 0) {
    $paraStart = $d.Paragraphs($targetIdx).Range.Start

    # Original run layout (character offsets relative to the paragraph start):
    #   0   "We "
    #   3   "regularly "
    #   13  "collaborate with other functions"
    #   45  " and our colleagues across RDQ"
    #   75  ". Here are some tasks "
    #   97  "where we work as a team"
    #   120 ". "
    #   122 <end>
    # Target layout: "We regularly collaborate" (bold, offsets 0-24) then
    # " with other functions" (not bold, offsets 24-45) as its own run,
    # with every run after that left completely untouched.

    $tailBoundaries = @(45, 75, 97, 120, 122)

    # Step 1: temporarily bold each of the untouched tail runs. This stops
    # the upcoming Find/Replace (step 2) from coalescing them into the
    # replacement run, since Word only merges *adjacent same-formatting*
    # runs on a text replace.
    for ($i = 0; $i -lt ($tailBoundaries.Count - 1); $i++) {
        $s = $paraStart + $tailBoundaries[$i]
        $e = $paraStart + $tailBoundaries[$i + 1]
        $d.Range($s, $e).Bold = 1
    }

    # Step 2: merge "We regularly collaborate" into a single bold run via
    # a Find & Replace (same text in both find/replace, only the
    # formatting changes) - this is what makes Word rebuild it as one run.
    $scope1 = $d.Range($paraStart, $paraStart + 24)
    $scope1.Find.ClearFormatting()
    $scope1.Find.Replacement.ClearFormatting()
    $scope1.Find.Replacement.Font.Bold = 1
    $scope1.Find.Execute("We regularly collaborate", $false, $false, $false, $false, $false, $true, 1, $false, "We regularly collaborate", 2) | Out-Null

    # Step 3: restore the tail runs back to their normal (non-bold) state,
    # one exact range at a time - plain property writes don't merge/split
    # anything beyond the range given.
    for ($i = 0; $i -lt ($tailBoundaries.Count - 1); $i++) {
        $s = $paraStart + $tailBoundaries[$i]
        $e = $paraStart + $tailBoundaries[$i + 1]
        $d.Range($s, $e).Bold = 0
    }

    # Step 4: the replace in step 2 leaves " with other functions" fused
    # to the (still-bold-marked-then-reset) tail; split it back out into
    # its own (non-bold) run with the same bold-toggle trick.
    $wrStart = $paraStart + 24
    $wrEnd = $paraStart + 45
    $wrRng = $d.Range($wrStart, $wrEnd)
    $wrRng.Bold = 1
    $wrRng.Bold = 0
}

# =========================================================================
# 2. Move the "_GoBack" bookmark from the "INCORRECT AUDIO FEEDBACK"
#    heading paragraph to the empty paragraph right after the
#    "We regularly collaborate ..." paragraph.
# =========================================================================

if ($targetIdx -gt 0) {
    $goBackPara = $d.Paragraphs($targetIdx + 1)
    $goBackRng = $goBackPara.Range
    $goBackRng.Collapse(1)

    # Bookmark names are unique, so adding a new "_GoBack" bookmark here
    # both removes the old one (next to "INCORRECT AUDIO FEEDBACK") and
    # creates the new one at this location.
    $d.Bookmarks.Add("_GoBack", $goBackRng)
}
